$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 14 (old row14 "total" -> row17, old row15 "footer" -> row18)
$ws.Rows("14:16").Insert()

# Copy formatting (incl. borders/merges/fonts) from the last data row (13) into the new rows
$ws.Range("A13:Q13").Copy()
$ws.Range("A14:Q16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 14: VOLTAREN 75MG/3ML 3 AMP.
$ws.Range("A14").Value = 8
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "VOLTAREN 75MG/3ML 3 AMP."
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "3:2"
$ws.Range("L14").NumberFormat = "@"
$ws.Range("L14").Value = "1"
$ws.Range("N14").NumberFormat = "@"
$ws.Range("N14").Value = "51.00"
$ws.Range("P14").NumberFormat = "@"
$ws.Range("P14").Value = "16.8300"
$ws.Range("Q14").NumberFormat = "@"
$ws.Range("Q14").Value = "0:1"

# Row 15: حبايه
$ws.Range("A15").Value = 9
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "حبايه"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "0:0"
$ws.Range("L15").NumberFormat = "@"
$ws.Range("L15").Value = "0"
$ws.Range("N15").NumberFormat = "@"
$ws.Range("N15").Value = "3.00"
$ws.Range("P15").NumberFormat = "@"
$ws.Range("P15").Value = "3.0000"
$ws.Range("Q15").NumberFormat = "@"
$ws.Range("Q15").Value = "1:0"

# Row 16: سرنجات 3 سم
$ws.Range("A16").Value = 10
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "سرنجات 3 سم"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "0:0"
$ws.Range("L16").NumberFormat = "@"
$ws.Range("L16").Value = "0"
$ws.Range("N16").NumberFormat = "@"
$ws.Range("N16").Value = "2.00"
$ws.Range("P16").NumberFormat = "@"
$ws.Range("P16").Value = "2.0000"
$ws.Range("Q16").NumberFormat = "@"
$ws.Range("Q16").Value = "1:0"

# Update the grand-total row (shifted from 14 to 17) with the new sum
$ws.Range("P17").Value = 481.83

# Update the footer timestamp (shifted from row 15 to row 18)
$ws.Range("A18").Value = "Tuesday, 10 June, 2025 9:41 AM"
